# corona_deaths_hist_israel.xlsx — "Add files via upload"
#
# The author appended 9 new daily-death records (rows 165-171 and 174 on
# the "10yr" sheet, with two blank spacer rows 172-173 kept in the same
# style as other blank spacer rows that already exist further up in the
# sheet). Everything else in the workbook (shared strings table, the
# lookup/summary tables in columns H:T, the "5yr" sheet, the two charts,
# AVERAGE/COUNTIF/COUNTIFS formulas, etc.) is a pure downstream
# consequence of those new rows and recalculates on its own — so the only
# thing this script does is add the rows, copying cell formatting from
# existing rows that already have the right look so new cells render
# exactly like their neighbours.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10yr")

# ---------------------------------------------------------------------
# 1. Prime formatting for the new rows by copying it from existing rows
#    that already look the way each new row should look. Column D
#    ("comments") is copied separately from A:C/E:F so that rows which
#    have no comment don't pick up a stray empty D cell:
#      - row 161  -> plain data row (date border-style age/sex/place,
#                    no comment, age not highlighted)
#      - row 147  -> same, but age is highlighted green/italic (s=23),
#                    used for the two rows with an approximate age
#      - row 162  -> blank spacer row (only the running index in col B)
#      - D126     -> plain (non-highlighted) "comments" cell style, used
#                    for the one new row that has a comment
# ---------------------------------------------------------------------

function Copy-RowFormat($srcRow, $dstRow) {
    $ws.Range("A${srcRow}:C${srcRow}").Copy() | Out-Null
    $ws.Range("A${dstRow}:C${dstRow}").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range("E${srcRow}:F${srcRow}").Copy() | Out-Null
    $ws.Range("E${dstRow}:F${dstRow}").PasteSpecial(-4122) | Out-Null
}

Copy-RowFormat 161 165
Copy-RowFormat 147 166
Copy-RowFormat 161 167
Copy-RowFormat 161 168
Copy-RowFormat 161 169

$ws.Range("D126").Copy() | Out-Null
$ws.Range("D169").PasteSpecial(-4122) | Out-Null

Copy-RowFormat 147 170
Copy-RowFormat 161 171

$ws.Range("B162").Copy() | Out-Null
$ws.Range("B172").PasteSpecial(-4122) | Out-Null
$ws.Range("B173").PasteSpecial(-4122) | Out-Null

Copy-RowFormat 161 174

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Fill in the actual new records.
#    Columns: A=date  B=running index (n-1)  C=age  D=comments
#             E=sex ("m"/"w")  F=place/hospital
# ---------------------------------------------------------------------

# 2020-04-18 (serial 43939)
$ws.Range("A165").Value = 43939
$ws.Range("B165").Value = 164
$ws.Range("C165").Value = 80
$ws.Range("D165").Value = "מחלות רקע שונות"
$ws.Range("E165").Value = "m"
$ws.Range("F165").Value = "הלל יפה"

$ws.Range("A166").Value = 43939
$ws.Range("B166").Value = 165
$ws.Range("C166").Value = 85
$ws.Range("E166").Value = "w"
$ws.Range("F166").Value = "פוריה"

# 2020-04-19 (serial 43940)
$ws.Range("A167").Value = 43940
$ws.Range("B167").Value = 166
$ws.Range("C167").Value = 96
$ws.Range("E167").Value = "m"
$ws.Range("F167").Value = "לניאדו"

$ws.Range("A168").Value = 43940
$ws.Range("B168").Value = 167
$ws.Range("C168").Value = 76
$ws.Range("E168").Value = "w"
$ws.Range("F168").Value = "הדסה עין כרם"

$ws.Range("A169").Value = 43940
$ws.Range("B169").Value = 168
$ws.Range("C169").Value = 29
$ws.Range("D169").Value = "מחלה ממארת סופנית"
$ws.Range("E169").Value = "w"
$ws.Range("F169").Value = "הדסה עין כרם"

$ws.Range("A170").Value = 43940
$ws.Range("B170").Value = 169
$ws.Range("C170").Value = 71
$ws.Range("E170").Value = "m"
$ws.Range("F170").Value = "פוריה"

$ws.Range("A171").Value = 43940
$ws.Range("B171").Value = 170
$ws.Range("C171").Value = 88
$ws.Range("E171").Value = "w"
$ws.Range("F171").Value = "אסותא אשדוד"

# Two blank spacer rows (as elsewhere in the sheet), only the running
# index is kept.
$ws.Range("B172").Value = 171
$ws.Range("B173").Value = 172

$ws.Range("A174").Value = 43940
$ws.Range("B174").Value = 173
$ws.Range("C174").Value = 87
$ws.Range("E174").Value = "w"
$ws.Range("F174").Value = "מעייני הישועה"

# ---------------------------------------------------------------------
# 3. Keep the selection/cursor roughly where the author left it.
# ---------------------------------------------------------------------
$ws.Range("C180").Select()

$wb.Save()
